$wb = $excel.ActiveWorkbook

# -------------------------------------------------------------------
# 1. accuracy_table sheet: restructure into separate "Statewide"/
#    "Clark County" tables plus a new pivoted "Accuracy Scores" block.
# -------------------------------------------------------------------
$acc = $wb.Worksheets.Item("accuracy_table")
$acc.Cells.Clear()

# --- table 1 (rows 2-7): statewide data, original scope column dropped
$acc.Range("A2").Value = "Target"
$acc.Range("B2").Value = "Number of Features"
$acc.Range("C2").Value = "Accuracy Score"
$acc.Range("D2").Value = "Rounded Accuracy Score"

$acc.Range("A3").Value = "Crash Severity"
$acc.Range("B3").Value = 1345
$acc.Range("C3").Value = 0.74000744456458101
$acc.Range("D3").Formula = "=ROUND(C3,2)"

$acc.Range("A4").Value = "Crash Severity"
$acc.Range("B4").Value = 773
$acc.Range("C4").Value = 0.73422904443696002
$acc.Range("D4").Formula = "=ROUND(C4,2)"

$acc.Range("A6").Value = "Injury Type"
$acc.Range("B6").Value = 1345
$acc.Range("C6").Value = 0.65668149671198395
$acc.Range("D6").Formula = "=ROUND(C6,2)"

$acc.Range("A7").Value = "Injury Type"
$acc.Range("B7").Value = 773
$acc.Range("C7").Value = 0.65538756048708702
$acc.Range("D7").Formula = "=ROUND(C7,2)"

# --- table 2 (rows 11-13): Clark County data
# New shared strings are introduced in the same order the original
# author typed them while finishing off the sheet.
$acc.Range("A18").Value = "Statewide Data"
$acc.Range("A11").Value = "Clark County Data"
$acc.Range("B18").Value = "All Columns"
$acc.Range("C18").Value = "Top 50%"
$acc.Range("A24").Value = "Clark County"
$acc.Range("A17").Value = "Accuracy Scores"
$acc.Range("D18").Value = "Top 16"

$acc.Range("A12").Value = "Crash Severity"
$acc.Range("B12").Value = 1232
$acc.Range("C12").Value = 0.74116122953715702
$acc.Range("D12").Formula = "=ROUND(C12,2)"

$acc.Range("A13").Value = "Injury Type"
$acc.Range("B13").Value = 1232
$acc.Range("C13").Value = 0.64609586621128201
$acc.Range("D13").Formula = "=ROUND(C13,2)"

# --- table 3 (rows 17-20): pivoted accuracy score summary
$acc.Range("A19").Value = "Injury Type"
$acc.Range("B19").Value = 0.66
$acc.Range("C19").Value = 0.66

$acc.Range("A20").Value = "Crash Severity"
$acc.Range("B20").Value = 0.74
$acc.Range("C20").Value = 0.73

# --- table 4 (rows 24-26): Clark County summary
$acc.Range("B24").Value = "All Columns"

$acc.Range("A25").Value = "Injury Type"
$acc.Range("B25").Value = 0.65

$acc.Range("A26").Value = "Crash Severity"
$acc.Range("B26").Value = 0.74

$acc.Columns.Item(1).ColumnWidth = 33
$acc.Columns.Item(2).ColumnWidth = 33
$acc.Columns.Item(3).ColumnWidth = 33
$acc.Columns.Item(4).ColumnWidth = 33

$acc.Range("B12").Select()

# -------------------------------------------------------------------
# 2. injury_type_scores sheet: sample size 1345 -> 1232 (Clark County)
# -------------------------------------------------------------------
$inj = $wb.Worksheets.Item("injury_type_scores")
$inj.Range("D5").ClearFormats()
$inj.Range("D5").Value = 1232
$inj.Range("D11").ClearFormats()
$inj.Range("D11").Value = 1232
$inj.Range("D5").Select()

# -------------------------------------------------------------------
# 3. crash_severity_scores sheet: sample size 1345 -> 1232
# -------------------------------------------------------------------
$sev = $wb.Worksheets.Item("crash_severity_scores")
$sev.Range("D6").Value = 1232
$sev.Range("D12").Value = 1232
$sev.Range("C19").Select()

# -------------------------------------------------------------------
# 4. Add a new "Sheet1" with feature-importance rankings
# -------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$s1 = $wb.Worksheets.Add($null, $lastSheet)
$s1.Name = "Sheet1"

$s1.Range("A2").Value = "latitude"
$s1.Range("A3").Value = "longitude"
$s1.Range("A4").Value = "v1_driver_age"
$s1.Range("A5").Value = "v2_driver_age"
$s1.Range("A6").Value = "crash_day"
$s1.Range("A7").Value = "crash_hour"
$s1.Range("A8").Value = "crash_month"
$s1.Range("A9").Value = "crash_day_of_week"
$s1.Range("A10").Value = "crash_year"
$s1.Range("A11").Value = "factors_roadway_1_UNKNOWN"
$s1.Range("A12").Value = "factors_roadway_1_DRY"
$s1.Range("A13").Value = "lighting_UNKNOWN"
$s1.Range("A14").Value = "total_vehicles"
$s1.Range("A15").Value = "lighting_DAYLIGHT"
$s1.Range("A16").Value = "hwy_factors_1_NONE"
$s1.Range("A17").Value = "v1_type_SEDAN_4DOOR"
$s1.Range("A18").Value = "hwy_factors_1_UNKNOWN"
$s1.Range("B1").Value = "injury_type"
$s1.Range("C1").Value = "crash_severity"

$s1.Range("B2").Value = 1
$s1.Range("C2").Value = 1
$s1.Range("B3").Value = 2
$s1.Range("C3").Value = 2
$s1.Range("B4").Value = 3
$s1.Range("C4").Value = 3
$s1.Range("B5").Value = 4
$s1.Range("C5").Value = 4
$s1.Range("B6").Value = 5
$s1.Range("C6").Value = 5
$s1.Range("B7").Value = 6
$s1.Range("C7").Value = 6
$s1.Range("B8").Value = 7
$s1.Range("C8").Value = 7
$s1.Range("B9").Value = 8
$s1.Range("C9").Value = 8
$s1.Range("B10").Value = 9
$s1.Range("C10").Value = 12
$s1.Range("B11").Value = 10
$s1.Range("C11").Value = 13
$s1.Range("B12").Value = 11
$s1.Range("C12").Value = 9
$s1.Range("B13").Value = 12
$s1.Range("C13").Value = 11
$s1.Range("B14").Value = 13
$s1.Range("C14").Value = 15
$s1.Range("B15").Value = 14
$s1.Range("C15").Value = 16
$s1.Range("B16").Value = 15
$s1.Range("C16").Value = 14
$s1.Range("B17").Value = 16
$s1.Range("C17").Value = 17
$s1.Range("B18").Value = 17
$s1.Range("C18").Value = 10

$s1.Columns.Item(1).ColumnWidth = 37
$s1.Columns.Item(2).ColumnWidth = 37

# re-apply the sort that produced this ranking (by injury_type rank, ascending)
$s1.Range("A2:C18").Sort($s1.Range("B2:B18"))

$s1.PageSetup.Orientation = 1

$s1.Range("E21").Select()
$s1.Activate()
